$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the computed prices in column D for rows 29-32
$ws.Range("D29").Value = 13023.612
$ws.Range("D30").Value = 15370.992
$ws.Range("D31").Value = 18326.952
$ws.Range("D32").Value = 21816.167
